# Apply the "disconnected_elements" diagnostic edit:
#   B1 = 0          (bold, centered/top, thin boxed border)
#   A2 = 0          (same style as B1)
#   B2 = "disconnected_elements"  (plain, shared string, no special style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- formatting for B1 -------------------------------------------------
$c1 = $ws.Range("B1")
$c1.Font.Bold = $true
$c1.HorizontalAlignment = -4108   # xlCenter
$c1.VerticalAlignment = -4160     # xlTop
$c1.Borders.LineStyle = 1         # xlContinuous
$c1.Borders.Weight = 2            # xlThin

# --- copy the same formatting onto A2 so both cells share one style ---
$c1.Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
